# Change the table style applied to the table on slide 16 (the "PLENARY -
# COMPLETE THE MISSING GAPS" slide) from the deck's custom "Table_0" style
# to a different (built-in) PowerPoint table style, identified by its
# brace-GUID StyleId.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shp = $s.Shapes.Item(3)
$tbl = $shp.Table
$tbl.ApplyStyle("{02827189-7717-4EF1-806F-009DBF2B6ED3}")
